$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "42.692.82"
$ws.Range("E2").Value = "  -0.63%  "

$ws.Range("D3").Value = "2.356.29"
$ws.Range("E3").Value = "  -1.30%  "

$ws.Range("E4").Value = "  -0.18%  "

Set-TextValue "D5" "318.08"
$ws.Range("E5").Value = "  -3.44%  "

Set-TextValue "D6" "109.48"
$ws.Range("E6").Value = "  +5.04%  "

Set-TextValue "D7" "0.634"
$ws.Range("E7").Value = "  -2.14%  "

$ws.Range("E8").Value = "  -0.05%  "

Set-TextValue "D9" "0.621"
$ws.Range("E9").Value = "  -4.71%  "

Set-TextValue "D10" "42.04"
$ws.Range("E10").Value = "  +0.77%  "

Set-TextValue "D11" "0.0928"
$ws.Range("E11").Value = "  -1.17%  "

Set-TextValue "D12" "8.63"
$ws.Range("E12").Value = "  -0.65%  "

$ws.Range("E13").Value = "  -4.50%  "

$ws.Range("E14").Value = "  +0.31%  "

Set-TextValue "D15" "16.09"
$ws.Range("E15").Value = "  -6.59%  "

$ws.Range("D16").Value = "2.705.60"
$ws.Range("E16").Value = "  -1.65%  "

$ws.Range("D17").Value = "2.506.68"
$ws.Range("E17").Value = "  +5.25%  "

$ws.Range("D18").Value = "42.661.50"
$ws.Range("E18").Value = "  -1.01%  "

Set-TextValue "D19" "7.74"
$ws.Range("E19").Value = "  -1.10%  "

Set-TextValue "D20" "0.0000107"
$ws.Range("E20").Value = "  -1.91%  "

$ws.Range("B21").Value = "PancakeSwap"
$ws.Range("C21").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D21" "3.78"
$ws.Range("E21").Value = "  +2.86%  "

$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D22" "76.02"
$ws.Range("E22").Value = "  -0.81%  "

Set-TextValue "D23" "256.70"
$ws.Range("E23").Value = "  -5.77%  "

Set-TextValue "D24" "2.34"
$ws.Range("E24").Value = "  -3.44%  "

Set-TextValue "D25" "9.47"
$ws.Range("E25").Value = "  -2.00%  "

$ws.Range("E26").Value = "  +0.06%  "

Set-TextValue "D27" "11.49"
$ws.Range("E27").Value = "  -2.20%  "

Set-TextValue "D28" "22.84"
$ws.Range("E28").Value = "  -0.68%  "

Set-TextValue "D29" "2.22"
$ws.Range("E29").Value = "  +1.60%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D30" "37.38"
$ws.Range("E30").Value = "  -0.48%  "

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D31" "173.39"
$ws.Range("E31").Value = "  -0.94%  "

Set-TextValue "D32" "0.0892"
$ws.Range("E32").Value = "  -4.15%  "

Set-TextValue "D33" "6.07"
$ws.Range("E33").Value = "  +3.07%  "

Set-TextValue "D34" "2.90"
$ws.Range("E34").Value = "  -8.39%  "

$ws.Range("E35").Value = "  +18.00%  "

$ws.Range("E36").Value = "  -1.94%  "

Set-TextValue "D37" "4.67"
$ws.Range("E37").Value = "  -5.76%  "

Set-TextValue "D38" "0.0364"
$ws.Range("E38").Value = "  -0.89%  "

Set-TextValue "D39" "3.94"
$ws.Range("E39").Value = "  -5.42%  "

Set-TextValue "D40" "2.71"
$ws.Range("E40").Value = "  -3.19%  "

Set-TextValue "D41" "0.239"
$ws.Range("E41").Value = "  +2.32%  "

Set-TextValue "D42" "1.49"
$ws.Range("E42").Value = "  -6.13%  "

Set-TextValue "D43" "70.81"
$ws.Range("E43").Value = "  +1.47%  "

$ws.Range("E44").Value = "  -0.15%  "

Set-TextValue "D45" "12.25"
$ws.Range("E45").Value = "  -1.16%  "

Set-TextValue "D46" "112.00"
$ws.Range("E46").Value = "  -7.98%  "

$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D47" "5.52"
$ws.Range("E47").Value = "  -0.27%  "

$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D48" "86.35"
$ws.Range("E48").Value = "  -4.54%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D49" "9.21"
$ws.Range("E49").Value = "  -1.67%  "

Set-TextValue "D50" "74.82"
$ws.Range("E50").Value = "  +1.46%  "

$ws.Range("E51").Value = "  -1.66%  "
